# Apply corrected IFRS figures (values appear to have been restated, likely
# from "thousand KRW" scale down to a different/corrected scale) and remove
# stray forecast rows (7-9) whose financial figures were erroneous, keeping
# only the label columns (A, B, C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New values for rows 2-6 (columns D..AJ) ---

$data = @{
    2 = @{
        D=2577; E=166; F=166; G=156; H=121; I=121; J=0; K=2060; L=1133; M=927;
        N=911; O=16; P=73; Q=101; R=-86; S=-170; T=65; U=36; V=284;
        W=6.44; X=4.7; Y=14.32; Z=5.88; AA=122.17; AB=956.12; AC=844; AD=3.01;
        AE=6666; AF=0.38; AG=20; AH=0.79; AI=2.26; AJ=14681617
    }
    3 = @{
        D=2630; E=161; F=161; G=191; H=152; I=151; J=0; K=2570; L=1467; M=1103;
        N=1086; O=17; P=79; Q=32; R=-162; S=80; T=122; U=-91; V=564;
        W=6.14; X=5.76; Y=15.14; Z=6.55; AA=132.94; AB=1095.06; AC=971; AD=4.67;
        AE=7316; AF=0.62; AG=30; AH=0.66; AI=2.95; AJ=15859354
    }
    4 = @{
        D=2643; E=115; F=115; G=114; H=92; I=93; J=-1; K=2525; L=1308; M=1216;
        N=1202; O=15; P=91; Q=188; R=-53; S=-17; T=108; U=80; V=541;
        W=4.37; X=3.49; Y=8.16; Z=3.62; AA=107.56; AB=1084.06; AC=548; AD=6.87;
        AE=7140; AF=0.53; AG=30; AH=0.8; AI=5.41; AJ=18275228
    }
    5 = @{
        D=2440; E=67; F=67; G=-2; H=-8; I=-8; J=0; K=2895; L=1711; M=1184;
        N=1170; O=14; P=92; Q=53; R=-369; S=384; T=417; U=-364; V=965;
        W=2.75; X=-0.32; Y=-0.66; Z=-0.29; AA=144.44; AB=1071.13; AC=-43; AD=-194.67;
        AE=7030; AF=1.18; AG=20; AH=0.24; AI=-42.61; AJ=18314054
    }
    6 = @{
        D=3030; E=-11; F=-11; G=-25; H=-40; I=-40; K=2972; L=1867; M=1105;
        N=1094; P=92; Q=-244; R=-29; S=165; T=61; U=-305; V=1173;
        W=-0.37; X=-1.33; Y=-3.51; Z=-1.38; AA=168.9; AB=988.5; AC=-217; AD=-30.66;
        AE=6572; AF=1.01; AG=10; AH=0.15; AI=-4.19; AJ=18314054
    }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}

# --- Rows 7-9: clear all financial figures (columns D..AJ), keep A, B, C ---

foreach ($row in 7..9) {
    $ws.Range("D$row`:AJ$row").ClearContents()
}
